# Complete map system redesign
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous sparse map layout (A1, B2, G14) before laying
# out the new, denser room/border map.
$ws.Cells.Clear()

# Row 1 - top room marker
$ws.Range("A1").Value = "room.1"

# Row 2 - border row across B:G
$ws.Range("B2:G2").Value = "b.1"

# Row 3 - borders, with two special intersection cells (b.1, b.3)
$ws.Range("B3").Value = "b.1"
$ws.Range("C3").Value = "b.1, b.3"
$ws.Range("D3").Value = "b.1, b.3"
$ws.Range("E3:G3").Value = "b.1"

# Row 4
$ws.Range("B4:G4").Value = "b.1"

# Row 5
$ws.Range("B5:G5").Value = "b.1"

# Row 6 - narrower row D:G (value written after row 10's special cell
# below, to match shared-string registration order)
$ws.Range("D6").Value = "b.1"
$ws.Range("F6:G6").Value = "b.1"

# Row 7
$ws.Range("D7:G7").Value = "b.1"

# Row 8 - widens out to D:I
$ws.Range("D8:I8").Value = "b.1"

# Row 9
$ws.Range("D9:I9").Value = "b.1"

# Row 10 - special cell (b.1, b.10) at H10
$ws.Range("D10:G10").Value = "b.1"
$ws.Range("H10").Value = "b.1, b.10"
$ws.Range("I10").Value = "b.1"

# Row 11
$ws.Range("D11:I11").Value = "b.1"

# Row 6 special cell (b.1, e.1), filled in after row 10 above
$ws.Range("E6").Value = "b.1, e.1"

# Row 12 - bottom room marker, shifted to column J
$ws.Range("J12").Value = "room.1"

# Restore the active selection to match the new layout
$ws.Range("E6").Select()
